$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.8
$ws.Range("H2").Value = 3.1
$ws.Range("I2").Value = 5.75
$ws.Range("J2").Value = 2.6
$ws.Range("L2").Value = 6.5
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 5
$ws.Range("X2").Value = 6.5
$ws.Range("AC2").Value = 5
$ws.Range("AD2").Value = 6.5
$ws.Range("AE2").Value = 26
$ws.Range("AI2").Value = 26
$ws.Range("AK2").Value = 67
$ws.Range("AL2").Value = 51
$ws.Range("AO2").Value = 11
$ws.Range("AP2").Value = 34
$ws.Range("AX2").Value = 34
$ws.Range("AZ2").Value = 151
$ws.Range("BA2").Value = 251
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 11
$ws.Range("BC3").Value = 151
$ws.Range("G4").Value = 1.73
$ws.Range("H4").Value = 3.3
$ws.Range("I4").Value = 5.75
$ws.Range("J4").Value = 2.4
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("Q4").Value = 2.35
$ws.Range("R4").Value = 1.57
$ws.Range("S4").Value = 1.5
$ws.Range("T4").Value = 2.5
$ws.Range("X4").Value = 7
$ws.Range("Y4").Value = 9
$ws.Range("Z4").Value = 13
$ws.Range("AA4").Value = 17
$ws.Range("AB4").Value = 34
$ws.Range("AC4").Value = 7
$ws.Range("AH4").Value = 11
$ws.Range("AI4").Value = 26
$ws.Range("AJ4").Value = 19
$ws.Range("AL4").Value = 51
$ws.Range("AN4").Value = 3.5
$ws.Range("AO4").Value = 9.5
$ws.Range("AQ4").Value = 34
$ws.Range("AT4").Value = 2.5
$ws.Range("AW4").Value = 6.5
$ws.Range("AZ4").Value = 126
$ws.Range("G6").Value = 2.77
$ws.Range("H6").Value = 2.87
$ws.Range("I6").Value = 2.65
$ws.Range("M6").Value = 1.02
$ws.Range("N6").Value = 7.5
$ws.Range("O6").Value = 1.36
$ws.Range("P6").Value = 2.7
$ws.Range("Q6").Value = 2.05
$ws.Range("R6").Value = 1.62
$ws.Range("U6").Value = 1.72
$ws.Range("V6").Value = 1.88
$ws.Range("X6").Value = 14.5
$ws.Range("AC6").Value = 7.7
$ws.Range("AD6").Value = 5.6
$ws.Range("AF6").Value = 65
$ws.Range("AH6").Value = 7.5
$ws.Range("AI6").Value = 13
$ws.Range("AJ6").Value = 9.75
$ws.Range("AL6").Value = 25
$ws.Range("AM6").Value = 35
